$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "78.969.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.177.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.96%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "632.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.228"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.177.89"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.582"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +32.85%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.44"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.765.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000225"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +15.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.36"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "78.921.35"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.179.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.42"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.95"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +27.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +12.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +13.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.30%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.50%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "76.42"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.90%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.92"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.23%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "518.94"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.66%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.136"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +20.16%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.88"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.123"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +15.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "164.35"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "191.21"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.817"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +14.17%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.58%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.50"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.40"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +12.25%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.622"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.21%  "
